$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(100, 8).Value = 11906465
$ws.Cells.Item(100, 9).Value = 15873955
$ws.Cells.Item(100, 11).Value = 15873955
$ws.Cells.Item(100, 13).Value = -15873414
$ws.Cells.Item(134, 8).Value = 44000
$ws.Cells.Item(134, 9).Value = 35000
$ws.Cells.Item(134, 10).Value = 50000
$ws.Cells.Item(134, 11).Value = 35000
$ws.Cells.Item(134, 12).Value = 50000
$ws.Cells.Item(134, 13).Value = -29930
$ws.Cells.Item(134, 14).Value = -60140
$ws.Cells.Item(138, 8).Value = 4402.143
$ws.Cells.Item(138, 9).Value = 1809.6
$ws.Cells.Item(138, 10).Value = 10883.5
$ws.Cells.Item(138, 11).Value = 5428.799999999999
$ws.Cells.Item(138, 12).Value = 32650.5
$ws.Cells.Item(138, 13).Value = -288.7999999999993
$ws.Cells.Item(138, 14).Value = -42930.5
$ws.Cells.Item(141, 8).Value = 4589.1055
$ws.Cells.Item(141, 9).Value = 4681.9414
$ws.Cells.Item(141, 10).Value = 3800
$ws.Cells.Item(141, 11).Value = 14045.8242
$ws.Cells.Item(141, 12).Value = 11400
$ws.Cells.Item(141, 13).Value = -8865.824199999999
$ws.Cells.Item(141, 14).Value = -21760

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 33476.188
$ws.Cells.Item(2, 9).Value = 46557.547
$ws.Cells.Item(2, 10).Value = 4697.2
$ws.Cells.Item(2, 11).Value = 46557.547
$ws.Cells.Item(2, 12).Value = 4697.2
$ws.Cells.Item(2, 13).Value = -46444.547
$ws.Cells.Item(2, 14).Value = -4923.2
$ws.Cells.Item(88, 8).Value = 5887.5
$ws.Cells.Item(88, 9).Value = 2000
$ws.Cells.Item(88, 10).Value = 6442.857
$ws.Cells.Item(88, 11).Value = 2000
$ws.Cells.Item(88, 12).Value = 6442.857
$ws.Cells.Item(88, 13).Value = -1594
$ws.Cells.Item(88, 14).Value = -7254.857
$ws.Cells.Item(91, 8).Value = 5887.5
$ws.Cells.Item(91, 9).Value = 2000
$ws.Cells.Item(91, 10).Value = 6442.857
$ws.Cells.Item(91, 11).Value = 2000
$ws.Cells.Item(91, 12).Value = 6442.857
$ws.Cells.Item(91, 13).Value = -596
$ws.Cells.Item(91, 14).Value = -9250.857
$ws.Cells.Item(102, 8).Value = 4012.0715
$ws.Cells.Item(102, 9).Value = 3418.7778
$ws.Cells.Item(102, 10).Value = 5080
$ws.Cells.Item(102, 11).Value = 3418.7778
$ws.Cells.Item(102, 12).Value = 5080
$ws.Cells.Item(102, 13).Value = -1796.7778
$ws.Cells.Item(102, 14).Value = -8324
$ws.Cells.Item(110, 8).Value = 192863.08
$ws.Cells.Item(110, 9).Value = 227839.19
$ws.Cells.Item(110, 10).Value = 494.5
$ws.Cells.Item(110, 11).Value = 227839.19
$ws.Cells.Item(110, 12).Value = 494.5
$ws.Cells.Item(110, 13).Value = -225794.19
$ws.Cells.Item(110, 14).Value = -4584.5
$ws.Cells.Item(116, 8).Value = 33476.188
$ws.Cells.Item(116, 9).Value = 46557.547
$ws.Cells.Item(116, 10).Value = 4697.2
$ws.Cells.Item(116, 11).Value = 46557.547
$ws.Cells.Item(116, 12).Value = 4697.2
$ws.Cells.Item(116, 13).Value = -44263.547
$ws.Cells.Item(116, 14).Value = -9285.200000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 33476.188
$ws.Cells.Item(3, 9).Value = 46557.547
$ws.Cells.Item(3, 10).Value = 4697.2
$ws.Cells.Item(3, 11).Value = 46557.547
$ws.Cells.Item(3, 12).Value = 4697.2
$ws.Cells.Item(3, 13).Value = -46443.547
$ws.Cells.Item(3, 14).Value = -4925.2
$ws.Cells.Item(99, 8).Value = 2429.625
$ws.Cells.Item(99, 9).Value = 1129.9231
$ws.Cells.Item(99, 10).Value = 3965.6365
$ws.Cells.Item(99, 11).Value = 1129.9231
$ws.Cells.Item(99, 12).Value = 3965.6365
$ws.Cells.Item(99, 13).Value = 368.0769
$ws.Cells.Item(99, 14).Value = -6961.636500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(80, 8).Value = 39600
$ws.Cells.Item(80, 10).Value = 39600
$ws.Cells.Item(80, 12).Value = 39600
$ws.Cells.Item(80, 14).Value = -41846
$ws.Cells.Item(83, 8).Value = 39600
$ws.Cells.Item(83, 10).Value = 39600
$ws.Cells.Item(83, 12).Value = 118800
$ws.Cells.Item(83, 14).Value = -130032
$ws.Cells.Item(94, 8).Value = 9956
$ws.Cells.Item(94, 9).Value = 0
$ws.Cells.Item(94, 10).Value = 9956
$ws.Cells.Item(94, 11).Value = 0
$ws.Cells.Item(94, 12).Value = 9956
$ws.Cells.Item(94, 13).ClearContents()
$ws.Cells.Item(94, 14).Value = -10858
$ws.Cells.Item(107, 8).Value = 37037450
$ws.Cells.Item(107, 9).Value = 52631970
$ws.Cells.Item(107, 10).Value = 467.5
$ws.Cells.Item(107, 11).Value = 52631970
$ws.Cells.Item(107, 12).Value = 467.5
$ws.Cells.Item(107, 13).Value = -52630050
$ws.Cells.Item(107, 14).Value = -4307.5
$ws.Cells.Item(112, 8).Value = 37700
$ws.Cells.Item(112, 10).Value = 37700
$ws.Cells.Item(112, 12).Value = 37700
$ws.Cells.Item(112, 14).Value = -40654

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(12, 8).Value = 157.2
$ws.Cells.Item(12, 10).Value = 260.58334
$ws.Cells.Item(12, 12).Value = 781.7500200000001
$ws.Cells.Item(12, 14).Value = -1127.75002
$ws.Cells.Item(34, 8).Value = 808.3684
$ws.Cells.Item(34, 9).Value = 103.2
$ws.Cells.Item(34, 10).Value = 1060.2142
$ws.Cells.Item(34, 11).Value = 309.6
$ws.Cells.Item(34, 12).Value = 3180.6426
$ws.Cells.Item(34, 13).Value = -225.6
$ws.Cells.Item(34, 14).Value = -3348.6426
$ws.Cells.Item(131, 8).Value = 5209264
$ws.Cells.Item(131, 9).Value = 646.2857
$ws.Cells.Item(131, 10).Value = 6098540
$ws.Cells.Item(131, 11).Value = 1938.8571
$ws.Cells.Item(131, 12).Value = 18295620
$ws.Cells.Item(131, 13).Value = 3101.1429
$ws.Cells.Item(131, 14).Value = -18305700

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(15, 8).Value = 14750
$ws.Cells.Item(15, 10).Value = 14750
$ws.Cells.Item(15, 12).Value = 14750
$ws.Cells.Item(15, 14).Value = -15326
$ws.Cells.Item(81, 8).Value = 14750
$ws.Cells.Item(81, 10).Value = 14750
$ws.Cells.Item(81, 12).Value = 14750
$ws.Cells.Item(81, 14).Value = -16746
$ws.Cells.Item(84, 8).Value = 14750
$ws.Cells.Item(84, 10).Value = 14750
$ws.Cells.Item(84, 12).Value = 44250
$ws.Cells.Item(84, 14).Value = -54234
$ws.Cells.Item(102, 8).Value = 2240.625
$ws.Cells.Item(102, 9).Value = 1753
$ws.Cells.Item(102, 11).Value = 1753
$ws.Cells.Item(102, 13).Value = -131
$ws.Cells.Item(126, 8).Value = 15153231
$ws.Cells.Item(126, 9).Value = 1985.1111
$ws.Cells.Item(126, 11).Value = 5955.3333
$ws.Cells.Item(126, 13).Value = -3485.3333

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 1889.4615
$ws.Cells.Item(7, 9).Value = 1976.5
$ws.Cells.Item(7, 10).Value = 1715.3846
$ws.Cells.Item(7, 11).Value = 1976.5
$ws.Cells.Item(7, 12).Value = 1715.3846
$ws.Cells.Item(7, 13).Value = -1864.5
$ws.Cells.Item(7, 14).Value = -1939.3846
$ws.Cells.Item(61, 8).Value = 2177.7
$ws.Cells.Item(61, 9).Value = 1371.0625
$ws.Cells.Item(61, 10).Value = 3099.5715
$ws.Cells.Item(61, 11).Value = 1371.0625
$ws.Cells.Item(61, 12).Value = 3099.5715
$ws.Cells.Item(61, 13).Value = -1169.0625
$ws.Cells.Item(61, 14).Value = -3503.5715
$ws.Cells.Item(100, 8).Value = 2259.2
$ws.Cells.Item(100, 9).Value = 1765.3334
$ws.Cells.Item(100, 10).Value = 3000
$ws.Cells.Item(100, 11).Value = 1765.3334
$ws.Cells.Item(100, 12).Value = 3000
$ws.Cells.Item(100, 13).Value = -1224.3334
$ws.Cells.Item(100, 14).Value = -4082
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()
$ws.Cells.Item(113, 8).Value = 2177.7
$ws.Cells.Item(113, 9).Value = 1371.0625
$ws.Cells.Item(113, 10).Value = 3099.5715
$ws.Cells.Item(113, 11).Value = 1371.0625
$ws.Cells.Item(113, 12).Value = 3099.5715
$ws.Cells.Item(113, 13).Value = 798.9375
$ws.Cells.Item(113, 14).Value = -7439.5715
$ws.Cells.Item(126, 8).Value = 1889.4615
$ws.Cells.Item(126, 9).Value = 1976.5
$ws.Cells.Item(126, 10).Value = 1715.3846
$ws.Cells.Item(126, 11).Value = 5929.5
$ws.Cells.Item(126, 12).Value = 5146.1538
$ws.Cells.Item(126, 13).Value = -3459.5
$ws.Cells.Item(126, 14).Value = -10086.1538
$ws.Cells.Item(132, 8).Value = 2305162.8
$ws.Cells.Item(132, 9).Value = 3759923.2
$ws.Cells.Item(132, 10).Value = 1791.9166
$ws.Cells.Item(132, 11).Value = 11279769.6
$ws.Cells.Item(132, 12).Value = 5375.7498
$ws.Cells.Item(132, 13).Value = -11277239.6
$ws.Cells.Item(132, 14).Value = -10435.7498
$ws.Cells.Item(136, 8).Value = 43524370
$ws.Cells.Item(136, 9).Value = 63013.062
$ws.Cells.Item(136, 10).Value = 142864600
$ws.Cells.Item(136, 11).Value = 189039.186
$ws.Cells.Item(136, 12).Value = 428593800
$ws.Cells.Item(136, 13).Value = -186489.186
$ws.Cells.Item(136, 14).Value = -428598900

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 32966.582
$ws.Cells.Item(136, 9).Value = 59143.883
$ws.Cells.Item(136, 11).Value = 177431.649
$ws.Cells.Item(136, 13).Value = -174881.649
